$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cell A1 = "Category", copying the header style (s="1")
# used by the rest of row 1 (e.g. B1), since there is no separate named style.
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Remove the header style from A2:A46 (revert to default/no style),
# while preserving their text content.
$ws.Range("A2:A46").ClearFormats()
